$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets("ALC")
$ws.Range("H9").Value = 469.09525
$ws.Range("I9").Value = 268.13333
$ws.Range("K9").Value = 268.13333
$ws.Range("M9").Value = -99.13333
$ws.Range("H28").Value = 847.1053000000001
$ws.Range("I28").Value = 843.73334
$ws.Range("J28").Value = 859.75
$ws.Range("K28").Value = 843.73334
$ws.Range("L28").Value = 859.75
$ws.Range("M28").Value = -358.73334
$ws.Range("N28").Value = -1829.75
$ws.Range("H43").Value = 2960.8333
$ws.Range("I43").Value = 2923.3333
$ws.Range("K43").Value = 2923.3333
$ws.Range("M43").Value = -2854.3333
$ws.Range("H62").Value = 46883204
$ws.Range("I62").Value = 19238750
$ws.Range("J62").Value = 166675840
$ws.Range("K62").Value = 19238750
$ws.Range("L62").Value = 166675840
$ws.Range("M62").Value = -19238126
$ws.Range("N62").Value = -166677088
$ws.Range("H65").Value = 46883204
$ws.Range("I65").Value = 19238750
$ws.Range("J65").Value = 166675840
$ws.Range("K65").Value = 96193750
$ws.Range("L65").Value = 833379200
$ws.Range("M65").Value = -96190630
$ws.Range("N65").Value = -833385440
$ws.Range("H69").Value = 9380.143
$ws.Range("J69").Value = 8391.75
$ws.Range("L69").Value = 25175.25
$ws.Range("N69").Value = -26923.25
$ws.Range("H70").Value = 7670.7744
$ws.Range("J70").Value = 14892.154
$ws.Range("L70").Value = 44676.462
$ws.Range("N70").Value = -45216.462
$ws.Range("H72").Value = 9380.143
$ws.Range("J72").Value = 8391.75
$ws.Range("L72").Value = 75525.75
$ws.Range("N72").Value = -84261.75
$ws.Range("H73").Value = 7670.7744
$ws.Range("J73").Value = 14892.154
$ws.Range("L73").Value = 44676.462
$ws.Range("N73").Value = -46548.462
$ws.Range("H80").Value = 658.44684
$ws.Range("J80").Value = 404.45456
$ws.Range("L80").Value = 1213.36368
$ws.Range("N80").Value = -3209.36368
$ws.Range("H83").Value = 658.44684
$ws.Range("J83").Value = 404.45456
$ws.Range("L83").Value = 3640.09104
$ws.Range("N83").Value = -13624.09104
$ws.Range("H86").Value = 1875.7693
$ws.Range("I86").Value = 2097
$ws.Range("J86").Value = 1737.5
$ws.Range("K86").Value = 2097
$ws.Range("L86").Value = 1737.5
$ws.Range("M86").Value = -974
$ws.Range("N86").Value = -3983.5
$ws.Range("H88").Value = 1534.6364
$ws.Range("I88").Value = 1538.8
$ws.Range("J88").Value = 1531.1666
$ws.Range("K88").Value = 1538.8
$ws.Range("L88").Value = 1531.1666
$ws.Range("M88").Value = -1132.8
$ws.Range("N88").Value = -2343.1666
$ws.Range("H89").Value = 1875.7693
$ws.Range("I89").Value = 2097
$ws.Range("J89").Value = 1737.5
$ws.Range("K89").Value = 10485
$ws.Range("L89").Value = 8687.5
$ws.Range("M89").Value = -4869
$ws.Range("N89").Value = -19919.5
$ws.Range("H91").Value = 1534.6364
$ws.Range("I91").Value = 1538.8
$ws.Range("J91").Value = 1531.1666
$ws.Range("K91").Value = 1538.8
$ws.Range("L91").Value = 1531.1666
$ws.Range("M91").Value = -134.8
$ws.Range("N91").Value = -4339.1666
$ws.Range("H98").Value = 5050149
$ws.Range("I98").Value = 5352528
$ws.Range("K98").Value = 5352528
$ws.Range("M98").Value = -5351030
$ws.Range("H106").Value = 2892.7856
$ws.Range("I106").Value = 2791.0833
$ws.Range("K106").Value = 2791.0833
$ws.Range("M106").Value = -2160.0833
$ws.Range("H111").Value = 3737.125
$ws.Range("I111").Value = 1816.1666
$ws.Range("K111").Value = 5448.4998
$ws.Range("M111").Value = -2381.4998
$ws.Range("H116").Value = 6526.1055
$ws.Range("J116").Value = 4978.6
$ws.Range("L116").Value = 4978.6
$ws.Range("N116").Value = -11862.6
$ws.Range("H122").Value = 5050149
$ws.Range("I122").Value = 5352528
$ws.Range("K122").Value = 16057584
$ws.Range("M122").Value = -16055134
$ws.Range("H125").Value = 8984.666999999999
$ws.Range("I125").Value = 2589.1428
$ws.Range("K125").Value = 23302.2852
$ws.Range("M125").Value = -20842.2852
$ws.Range("H129").Value = 3376.4092
$ws.Range("I129").Value = 2431.6667
$ws.Range("K129").Value = 7295.000100000001
$ws.Range("M129").Value = -2295.000100000001
$ws.Range("H132").Value = 4445.25
$ws.Range("I132").Value = 3901.6667
$ws.Range("K132").Value = 11705.0001
$ws.Range("M132").Value = -9175.000100000001
$ws.Range("H137").Value = 13883.5
$ws.Range("J137").Value = 4734
$ws.Range("L137").Value = 14202
$ws.Range("N137").Value = -19302
$ws.Range("H138").Value = 4604.64
$ws.Range("J138").Value = 4868.282
$ws.Range("L138").Value = 14604.846
$ws.Range("N138").Value = -24884.846

$ws = $wb.Sheets("ARM")
$ws.Range("H2").Value = 7000.607
$ws.Range("I2").Value = 7394.278
$ws.Range("K2").Value = 7394.278
$ws.Range("M2").Value = -7281.278
$ws.Range("H5").Value = 221.75
$ws.Range("I5").Value = 229.83333
$ws.Range("J5").Value = 197.5
$ws.Range("K5").Value = 229.83333
$ws.Range("L5").Value = 197.5
$ws.Range("M5").Value = -117.83333
$ws.Range("N5").Value = -421.5
$ws.Range("H32").Value = 10408.921
$ws.Range("I32").Value = 10037.054
$ws.Range("J32").Value = 38299
$ws.Range("K32").Value = 10037.054
$ws.Range("L32").Value = 38299
$ws.Range("M32").Value = -9750.054
$ws.Range("N32").Value = -38873
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H45").Value = 2638.913
$ws.Range("I45").Value = 1235
$ws.Range("K45").Value = 1235
$ws.Range("M45").Value = -858
$ws.Range("H63").Value = 5335
$ws.Range("I63").Value = 5502.5
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 5502.5
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -4816.5
$ws.Range("N63").Value = -6372
$ws.Range("H66").Value = 5335
$ws.Range("I66").Value = 5502.5
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 27512.5
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -24080.5
$ws.Range("N66").Value = -31864
$ws.Range("H74").Value = 4822.0317
$ws.Range("I74").Value = 4810.1187
$ws.Range("J74").Value = 4997.75
$ws.Range("K74").Value = 4810.1187
$ws.Range("L74").Value = 4997.75
$ws.Range("M74").Value = -3936.1187
$ws.Range("N74").Value = -6745.75
$ws.Range("H77").Value = 4822.0317
$ws.Range("I77").Value = 4810.1187
$ws.Range("J77").Value = 4997.75
$ws.Range("K77").Value = 24050.5935
$ws.Range("L77").Value = 24988.75
$ws.Range("M77").Value = -19682.5935
$ws.Range("N77").Value = -33724.75
$ws.Range("H80").Value = 39999.5
$ws.Range("I80").Value = 39999.5
$ws.Range("K80").Value = 39999.5
$ws.Range("M80").Value = -39001.5
$ws.Range("H83").Value = 39999.5
$ws.Range("I83").Value = 39999.5
$ws.Range("K83").Value = 119998.5
$ws.Range("M83").Value = -115006.5
$ws.Range("H97").Value = 1212.25
$ws.Range("I97").Value = 722.8461
$ws.Range("J97").Value = 3333
$ws.Range("K97").Value = 722.8461
$ws.Range("L97").Value = 3333
$ws.Range("M97").Value = -226.8461
$ws.Range("N97").Value = -4325
$ws.Range("H102").Value = 3258.389
$ws.Range("J102").Value = 3734.7144
$ws.Range("L102").Value = 3734.7144
$ws.Range("N102").Value = -6978.7144
$ws.Range("H110").Value = 27278.666
$ws.Range("I110").Value = 38124.75
$ws.Range("J110").Value = 18601.8
$ws.Range("K110").Value = 38124.75
$ws.Range("L110").Value = 18601.8
$ws.Range("M110").Value = -36079.75
$ws.Range("N110").Value = -22691.8
$ws.Range("H116").Value = 7000.607
$ws.Range("I116").Value = 7394.278
$ws.Range("K116").Value = 7394.278
$ws.Range("M116").Value = -5100.278
$ws.Range("H132").Value = 4784.5
$ws.Range("I132").Value = 3723.0715
$ws.Range("K132").Value = 11169.2145
$ws.Range("M132").Value = -8639.2145

$ws = $wb.Sheets("BSM")
$ws.Range("H3").Value = 7000.607
$ws.Range("I3").Value = 7394.278
$ws.Range("K3").Value = 7394.278
$ws.Range("M3").Value = -7280.278
$ws.Range("H4").Value = 221.75
$ws.Range("I4").Value = 229.83333
$ws.Range("J4").Value = 197.5
$ws.Range("K4").Value = 229.83333
$ws.Range("L4").Value = 197.5
$ws.Range("M4").Value = -114.83333
$ws.Range("N4").Value = -427.5
$ws.Range("H20").Value = 2405.027
$ws.Range("I20").Value = 1406.6666
$ws.Range("K20").Value = 1406.6666
$ws.Range("M20").Value = -1159.6666
$ws.Range("H36").Value = 2824.75
$ws.Range("I36").Value = 2824.75
$ws.Range("K36").Value = 2824.75
$ws.Range("M36").Value = -2290.75
$ws.Range("H40").Value = 39999
$ws.Range("J40").Value = 39999
$ws.Range("L40").Value = 39999
$ws.Range("N40").Value = -40529
$ws.Range("H76").Value = 40157
$ws.Range("I76").Value = 30000
$ws.Range("J76").Value = 50314
$ws.Range("K76").Value = 30000
$ws.Range("L76").Value = 50314
$ws.Range("M76").Value = -29685
$ws.Range("N76").Value = -50944
$ws.Range("H79").Value = 40157
$ws.Range("I79").Value = 30000
$ws.Range("J79").Value = 50314
$ws.Range("K79").Value = 30000
$ws.Range("L79").Value = 50314
$ws.Range("M79").Value = -28908
$ws.Range("N79").Value = -52498
$ws.Range("H86").Value = 2600.8333
$ws.Range("I86").Value = 1921.2
$ws.Range("J86").Value = 5999
$ws.Range("K86").Value = 1921.2
$ws.Range("L86").Value = 5999
$ws.Range("M86").Value = -798.2
$ws.Range("N86").Value = -8245
$ws.Range("H89").Value = 2600.8333
$ws.Range("I89").Value = 1921.2
$ws.Range("J89").Value = 5999
$ws.Range("K89").Value = 9606
$ws.Range("L89").Value = 29995
$ws.Range("M89").Value = -3990
$ws.Range("N89").Value = -41227
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").ClearContents()
$ws.Range("H96").Value = 166679680
$ws.Range("I96").Value = 200003800
$ws.Range("J96").Value = 59000
$ws.Range("K96").Value = 200003800
$ws.Range("L96").Value = 59000
$ws.Range("M96").Value = -200001054
$ws.Range("N96").Value = -64492
$ws.Range("H99").Value = 8728.200000000001
$ws.Range("I99").Value = 17626.143
$ws.Range("J99").Value = 3937
$ws.Range("K99").Value = 17626.143
$ws.Range("L99").Value = 3937
$ws.Range("M99").Value = -16128.143
$ws.Range("N99").Value = -6933
$ws.Range("H134").Value = 2200.8442
$ws.Range("I134").Value = 1945.7142
$ws.Range("J134").Value = 2413.4524
$ws.Range("K134").Value = 5837.142599999999
$ws.Range("L134").Value = 7240.3572
$ws.Range("M134").Value = -3302.142599999999
$ws.Range("N134").Value = -12310.3572
$ws.Range("H140").Value = 78000
$ws.Range("J140").Value = 78000
$ws.Range("L140").Value = 78000
$ws.Range("N140").Value = -88360

$ws = $wb.Sheets("CRP")
$ws.Range("H9").Value = 190137
$ws.Range("J9").Value = 190137
$ws.Range("L9").Value = 190137
$ws.Range("N9").Value = -190473
$ws.Range("H16").Value = 3724.5186
$ws.Range("I16").Value = 3618.25
$ws.Range("K16").Value = 3618.25
$ws.Range("M16").Value = -3331.25
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H62").Value = 20568.572
$ws.Range("J62").Value = 2999
$ws.Range("L62").Value = 2999
$ws.Range("N62").Value = -4247
$ws.Range("H65").Value = 20568.572
$ws.Range("J65").Value = 2999
$ws.Range("L65").Value = 14995
$ws.Range("N65").Value = -21235
$ws.Range("H92").Value = 44559.96
$ws.Range("J92").Value = 44559.96
$ws.Range("L92").Value = 44559.96
$ws.Range("N92").Value = -49551.96
$ws.Range("H107").Value = 5176.7417
$ws.Range("I107").Value = 715.5
$ws.Range("K107").Value = 715.5
$ws.Range("M107").Value = 1204.5
$ws.Range("H113").Value = 3724.5186
$ws.Range("I113").Value = 3618.25
$ws.Range("K113").Value = 3618.25
$ws.Range("M113").Value = -1448.25
$ws.Range("H122").Value = 3014.6553
$ws.Range("I122").Value = 2892.2856
$ws.Range("J122").Value = 3335.875
$ws.Range("K122").Value = 8676.856800000001
$ws.Range("L122").Value = 10007.625
$ws.Range("M122").Value = -6226.856800000001
$ws.Range("N122").Value = -14907.625
$ws.Range("H132").Value = 4311.2104
$ws.Range("I132").Value = 3620.625
$ws.Range("K132").Value = 10861.875
$ws.Range("M132").Value = -8331.875
$ws.Range("H134").Value = 5276.212
$ws.Range("I134").Value = 5383.3105
$ws.Range("J134").Value = 4499.75
$ws.Range("K134").Value = 16149.9315
$ws.Range("L134").Value = 13499.25
$ws.Range("M134").Value = -13614.9315
$ws.Range("N134").Value = -18569.25

$ws = $wb.Sheets("CUL")
$ws.Range("H2").Value = 6250412
$ws.Range("I2").Value = 553.1667
$ws.Range("J2").Value = 10000327
$ws.Range("K2").Value = 3319.0002
$ws.Range("L2").Value = 60001962
$ws.Range("M2").Value = -3206.0002
$ws.Range("N2").Value = -60002188
$ws.Range("H7").Value = 201552
$ws.Range("I7").Value = 500125
$ws.Range("J7").Value = 2503.3333
$ws.Range("K7").Value = 1500375
$ws.Range("L7").Value = 7509.999899999999
$ws.Range("M7").Value = -1500263
$ws.Range("N7").Value = -7733.999899999999
$ws.Range("H17").Value = 130.25
$ws.Range("J17").Value = 300
$ws.Range("L17").Value = 900
$ws.Range("N17").Value = -1238
$ws.Range("H22").Value = 1594.6
$ws.Range("I22").Value = 980
$ws.Range("K22").Value = 2940
$ws.Range("M22").Value = -2771
$ws.Range("H23").Value = 191.83333
$ws.Range("J23").Value = 274
$ws.Range("L23").Value = 822
$ws.Range("N23").Value = -1292
$ws.Range("H26").Value = 303.33334
$ws.Range("J26").Value = 299.85715
$ws.Range("L26").Value = 899.5714499999999
$ws.Range("N26").Value = -1475.57145
$ws.Range("H27").Value = 1594.6
$ws.Range("I27").Value = 980
$ws.Range("K27").Value = 2940
$ws.Range("M27").Value = -2838
$ws.Range("H40").Value = 55
$ws.Range("I40").Value = 67.14286
$ws.Range("J40").Value = 33.75
$ws.Range("K40").Value = 268.57144
$ws.Range("L40").Value = 135
$ws.Range("M40").Value = -199.57144
$ws.Range("N40").Value = -273
$ws.Range("H55").Value = 6710
$ws.Range("J55").Value = 9850
$ws.Range("L55").Value = 29550
$ws.Range("N55").Value = -29904
$ws.Range("H93").Value = 49999.5
$ws.Range("I93").Value = 89999
$ws.Range("J93").Value = 10000
$ws.Range("K93").Value = 269997
$ws.Range("L93").Value = 30000
$ws.Range("M93").Value = -268125
$ws.Range("N93").Value = -33744
$ws.Range("H131").Value = 13159476
$ws.Range("I131").Value = 90909670
$ws.Range("J131").Value = 1750.9846
$ws.Range("K131").Value = 272729010
$ws.Range("L131").Value = 5252.9538
$ws.Range("M131").Value = -272723970
$ws.Range("N131").Value = -15332.9538
$ws.Range("H132").Value = 4282.3022
$ws.Range("I132").Value = 6062.04
$ws.Range("K132").Value = 54558.36
$ws.Range("M132").Value = -52028.36

$ws = $wb.Sheets("GSM")
$ws.Range("H70").Value = 76929440
$ws.Range("I70").Value = 5470.143
$ws.Range("J70").Value = 166674060
$ws.Range("K70").Value = 5470.143
$ws.Range("L70").Value = 166674060
$ws.Range("M70").Value = -5200.143
$ws.Range("N70").Value = -166674600
$ws.Range("H73").Value = 76929440
$ws.Range("I73").Value = 5470.143
$ws.Range("J73").Value = 166674060
$ws.Range("K73").Value = 5470.143
$ws.Range("L73").Value = 166674060
$ws.Range("M73").Value = -4534.143
$ws.Range("N73").Value = -166675932
$ws.Range("H80").Value = 3133.4285
$ws.Range("I80").Value = 3846.8
$ws.Range("J80").Value = 1350
$ws.Range("K80").Value = 3846.8
$ws.Range("L80").Value = 1350
$ws.Range("M80").Value = -2848.8
$ws.Range("N80").Value = -3346
$ws.Range("H83").Value = 3133.4285
$ws.Range("I83").Value = 3846.8
$ws.Range("J83").Value = 1350
$ws.Range("K83").Value = 19234
$ws.Range("L83").Value = 6750
$ws.Range("M83").Value = -14242
$ws.Range("N83").Value = -16734
$ws.Range("H97").Value = 1808.5454
$ws.Range("I97").Value = 1489.4
$ws.Range("J97").Value = 5000
$ws.Range("K97").Value = 1489.4
$ws.Range("L97").Value = 5000
$ws.Range("M97").Value = -993.4000000000001
$ws.Range("N97").Value = -5992
$ws.Range("H113").Value = 8609.958000000001
$ws.Range("I113").Value = 9375.888999999999
$ws.Range("K113").Value = 9375.888999999999
$ws.Range("M113").Value = -7205.888999999999
$ws.Range("H122").Value = 4067.625
$ws.Range("I122").Value = 3108.9
$ws.Range("K122").Value = 9326.700000000001
$ws.Range("M122").Value = -6876.700000000001
$ws.Range("H126").Value = 12222.25
$ws.Range("I126").Value = 14870.6
$ws.Range("K126").Value = 44611.8
$ws.Range("M126").Value = -42141.8

$ws = $wb.Sheets("LTW")
$ws.Range("H16").Value = 22730816
$ws.Range("J16").Value = 10178.6
$ws.Range("L16").Value = 10178.6
$ws.Range("N16").Value = -10518.6
$ws.Range("H22").Value = 2020.5454
$ws.Range("J22").Value = 2385.3333
$ws.Range("L22").Value = 2385.3333
$ws.Range("N22").Value = -2975.3333
$ws.Range("H27").Value = 2020.5454
$ws.Range("J27").Value = 2385.3333
$ws.Range("L27").Value = 2385.3333
$ws.Range("N27").Value = -2599.3333
$ws.Range("H34").Value = 4074.1428
$ws.Range("I34").Value = 6498.3335
$ws.Range("J34").Value = 2256
$ws.Range("K34").Value = 6498.3335
$ws.Range("L34").Value = 2256
$ws.Range("M34").Value = -6326.3335
$ws.Range("N34").Value = -2600
$ws.Range("H42").Value = 19666
$ws.Range("I42").Value = 21999
$ws.Range("J42").Value = 15000
$ws.Range("K42").Value = 21999
$ws.Range("L42").Value = 15000
$ws.Range("M42").Value = -21436
$ws.Range("N42").Value = -16126
$ws.Range("H46").Value = 2612.2
$ws.Range("I46").Value = 1398.6923
$ws.Range("J46").Value = 10500
$ws.Range("K46").Value = 1398.6923
$ws.Range("L46").Value = 10500
$ws.Range("M46").Value = -1210.6923
$ws.Range("N46").Value = -10876
$ws.Range("H49").Value = 19666
$ws.Range("I49").Value = 21999
$ws.Range("J49").Value = 15000
$ws.Range("K49").Value = 21999
$ws.Range("L49").Value = 15000
$ws.Range("M49").Value = -21852
$ws.Range("N49").Value = -15294
$ws.Range("H61").Value = 1675.6666
$ws.Range("I61").Value = 1771.8182
$ws.Range("J61").Value = 1252.6
$ws.Range("K61").Value = 1771.8182
$ws.Range("L61").Value = 1252.6
$ws.Range("M61").Value = -1569.8182
$ws.Range("N61").Value = -1656.6
$ws.Range("H82").Value = 2505.2
$ws.Range("J82").Value = 2138.25
$ws.Range("L82").Value = 2138.25
$ws.Range("N82").Value = -2860.25
$ws.Range("H85").Value = 2505.2
$ws.Range("J85").Value = 2138.25
$ws.Range("L85").Value = 2138.25
$ws.Range("N85").Value = -4634.25
$ws.Range("H93").Value = 876.2222
$ws.Range("I93").Value = 860.875
$ws.Range("J93").Value = 999
$ws.Range("K93").Value = 860.875
$ws.Range("L93").Value = 999
$ws.Range("M93").Value = 387.125
$ws.Range("N93").Value = -3495
$ws.Range("H113").Value = 1675.6666
$ws.Range("I113").Value = 1771.8182
$ws.Range("J113").Value = 1252.6
$ws.Range("K113").Value = 1771.8182
$ws.Range("L113").Value = 1252.6
$ws.Range("M113").Value = 398.1818000000001
$ws.Range("N113").Value = -5592.6
$ws.Range("H122").Value = 7270.4287
$ws.Range("I122").Value = 5754.5557
$ws.Range("K122").Value = 17263.6671
$ws.Range("M122").Value = -14813.6671
$ws.Range("H123").Value = 74714
$ws.Range("J123").Value = 74714
$ws.Range("L123").Value = 74714
$ws.Range("N123").Value = -84514
$ws.Range("H136").Value = 2632.5625
$ws.Range("I136").Value = 2327.5417
$ws.Range("J136").Value = 3547.625
$ws.Range("K136").Value = 6982.625100000001
$ws.Range("L136").Value = 10642.875
$ws.Range("M136").Value = -4432.625100000001
$ws.Range("N136").Value = -15742.875

$ws = $wb.Sheets("WVR")
$ws.Range("H15").Value = 14499
$ws.Range("J15").Value = 14499
$ws.Range("L15").Value = 14499
$ws.Range("N15").Value = -15075
$ws.Range("H25").Value = 17500
$ws.Range("I25").Value = 15000
$ws.Range("K25").Value = 15000
$ws.Range("M25").Value = -14707
$ws.Range("H43").Value = 33979.6
$ws.Range("I43").Value = 34974.5
$ws.Range("K43").Value = 34974.5
$ws.Range("M43").Value = -34825.5
$ws.Range("H54").Value = 27249.5
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H62").Value = 15636055
$ws.Range("J62").Value = 19241914
$ws.Range("L62").Value = 19241914
$ws.Range("N62").Value = -19243162
$ws.Range("H65").Value = 15636055
$ws.Range("J65").Value = 19241914
$ws.Range("L65").Value = 96209570
$ws.Range("N65").Value = -96215810
$ws.Range("H112").Value = 28890
$ws.Range("J112").Value = 30500
$ws.Range("L112").Value = 30500
$ws.Range("N112").Value = -33454
$ws.Range("H122").Value = 2002.2
$ws.Range("I122").Value = 2002.75
$ws.Range("K122").Value = 6008.25
$ws.Range("M122").Value = -3558.25
$ws.Range("H126").Value = 4215.421
$ws.Range("I126").Value = 2463.9285
$ws.Range("J126").Value = 9119.6
$ws.Range("K126").Value = 7391.7855
$ws.Range("L126").Value = 27358.8
$ws.Range("M126").Value = -4921.7855
$ws.Range("N126").Value = -32298.8
$ws.Range("H132").Value = 16395432
$ws.Range("I132").Value = 24391948
$ws.Range("K132").Value = 73175844
$ws.Range("M132").Value = -73173314
$ws.Range("H136").Value = 2468.0505
$ws.Range("I136").Value = 1640.7742
$ws.Range("J136").Value = 3854.2974
$ws.Range("K136").Value = 4922.3226
$ws.Range("L136").Value = 11562.8922
$ws.Range("M136").Value = -2372.3226
$ws.Range("N136").Value = -16662.8922
